# Automatic update of files.
# Bump the "Förändrad" (Changed) date in column C from 2023-09-20 (45189)
# to 2023-09-21 (45190) for every data row (rows 2-181).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C181").Value = 45190
